# Adds an "xfile" attribute row to the "attributes" sheet (describing the
# new xfile column on org_molgenis_test_TypeTest) and adds the matching
# "xfile" header column to the org_molgenis_test_TypeTest sheet.

$wb = $excel.ActiveWorkbook

# --- org_molgenis_test_TypeTest: add new "xfile" header column (AS1) ---
$wsType = $wb.Worksheets.Item("org_molgenis_test_TypeTest")
$wsType.Range("AS1").Value = "xfile"
$wsType.Range("AT26").Select()

# --- attributes: insert a new row describing the xfile attribute ---
$wsAttr = $wb.Worksheets.Item("attributes")
$wsAttr.Activate()
$wsAttr.Rows.Item(49).Insert()

$wsAttr.Range("A49").Value = "xfile"
$wsAttr.Range("B49").Value = "org_molgenis_test_TypeTest"
$wsAttr.Range("C49").Value = "file"
$wsAttr.Range("E49").Value = $false
$wsAttr.Range("F49").Value = $true

$wsAttr.Range("G48").Select()
